$d = $word.ActiveDocument

$replacements = @(
    @("487÷9=", "299÷6="),
    @("693÷6=", "459÷6="),
    @("149÷4=", "513÷7="),
    @("131÷6=", "657÷8="),
    @("176÷3=", "327÷5="),
    @("483÷3=", "777÷3="),
    @("862÷7=", "178÷3="),
    @("401÷3=", "762÷4="),
    @("603÷9=", "624÷5="),
    @("735÷2=", "278÷8="),
    @("669÷2=", "365÷4="),
    @("555÷7=", "332÷3="),
    @("581÷5=", "644÷7="),
    @("160÷9=", "636÷9="),
    @("658÷3=", "775÷6="),
    @("152÷4=", "597÷5="),
    @("992÷9=", "557÷2="),
    @("939÷9=", "970÷6="),
    @("973÷9=", "675÷2="),
    @("116÷6=", "894÷8="),
    @("108÷8=", "491÷7="),
    @("765÷7=", "777÷5="),
    @("955÷2=", "744÷7="),
    @("205÷2=", "891÷8="),
    @("480÷6=", "726÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
